# The presentation ships with two theme parts:
#   - the theme applied to the Slide Master ("Integral" / "Red Violet" palette)
#   - the theme applied to the Notes Master (the stock "Office Theme" palette)
#
# The authored edit swaps which palette each master uses: the Slide Master
# (and therefore the whole deck's look) is switched over to the plain
# default "Office Theme" colors, while the Notes Master keeps the palette
# that used to belong to the slides ("Integral" / "Red Violet").
#
# Re-create that by pushing the default Office theme's 12 color-scheme
# slots onto the presentation's (Slide Master) theme color scheme, in the
# standard dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order.

function ConvertTo-OfficeRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - stock "Office Theme" colors
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47",
    "0563C1", "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-OfficeRgb $officeThemeColors[$i - 1]
}
